$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "0.995", "214.29") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.783.94"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "1.629.58"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("D4").Value = "0.995"
$ws.Range("E4").Value = "  -0.86%  "

$ws.Range("D5").Value = "214.29"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.854.60"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.623.48"
$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "0.0₃0762"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "62.84"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "25.783.08"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").Value = "191.38"
$ws.Range("E21").Value = "  -1.13%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("E25").Value = "  +1.86%  "

$ws.Range("D26").Value = "141.92"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("E27").Value = "  +2.73%  "

$ws.Range("D28").Value = "6.85"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +2.55%  "

$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").Value = "1.145.85"
$ws.Range("E37").Value = "  +3.87%  "

$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").Value = "0.995"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("D44").Value = "100.72"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("E45").Value = "  +1.32%  "

$ws.Range("D46").Value = "1.765.38"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").Value = "55.34"
$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.46"
$ws.Range("E50").Value = "  +6.28%  "

$ws.Range("E51").Value = "  -0.54%  "

# Restore the default (Normal) style for column D so only the
# values change and no residual number-format styling remains.
$ws.Range("D2:D51").Style = "Normal"
